$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '29.424.64'; Numeric = $false }
    @{ Cell = 'E2'; Value = '  +0.62%  '; Numeric = $false }
    @{ Cell = 'D3'; Value = '1.872.13'; Numeric = $false }
    @{ Cell = 'E3'; Value = '  +0.91%  '; Numeric = $false }
    @{ Cell = 'D4'; Value = '1.001'; Numeric = $true }
    @{ Cell = 'E4'; Value = '  +0.00%  '; Numeric = $false }
    @{ Cell = 'D5'; Value = '0.7115'; Numeric = $true }
    @{ Cell = 'E5'; Value = '  +1.40%  '; Numeric = $false }
    @{ Cell = 'D6'; Value = '242.26'; Numeric = $true }
    @{ Cell = 'E6'; Value = '  +1.87%  '; Numeric = $false }
    @{ Cell = 'D7'; Value = '1.001'; Numeric = $true }
    @{ Cell = 'E7'; Value = '  +0.00%  '; Numeric = $false }
    @{ Cell = 'D8'; Value = '0.07865'; Numeric = $true }
    @{ Cell = 'E8'; Value = '  -2.05%  '; Numeric = $false }
    @{ Cell = 'D9'; Value = '0.3109'; Numeric = $true }
    @{ Cell = 'E9'; Value = '  +2.88%  '; Numeric = $false }
    @{ Cell = 'D10'; Value = '24.98'; Numeric = $true }
    @{ Cell = 'E10'; Value = '  +6.20%  '; Numeric = $false }
    @{ Cell = 'D11'; Value = '0.08250'; Numeric = $true }
    @{ Cell = 'E11'; Value = '  +0.89%  '; Numeric = $false }
    @{ Cell = 'D12'; Value = '1.883.10'; Numeric = $false }
    @{ Cell = 'E12'; Value = '  +1.26%  '; Numeric = $false }
    @{ Cell = 'B13'; Value = 'Polkadot'; Numeric = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; Numeric = $false }
    @{ Cell = 'D13'; Value = '5.291'; Numeric = $true }
    @{ Cell = 'E13'; Value = '  +1.71%  '; Numeric = $false }
    @{ Cell = 'B14'; Value = 'Polygon'; Numeric = $false }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; Numeric = $false }
    @{ Cell = 'D14'; Value = '0.7255'; Numeric = $true }
    @{ Cell = 'E14'; Value = '  +2.68%  '; Numeric = $false }
    @{ Cell = 'D15'; Value = '90.77'; Numeric = $true }
    @{ Cell = 'E15'; Value = '  +1.12%  '; Numeric = $false }
    @{ Cell = 'D16'; Value = '29.450.98'; Numeric = $false }
    @{ Cell = 'E16'; Value = '  +0.61%  '; Numeric = $false }
    @{ Cell = 'D17'; Value = '5.924'; Numeric = $true }
    @{ Cell = 'E17'; Value = '  +1.57%  '; Numeric = $false }
    @{ Cell = 'D18'; Value = '247.71'; Numeric = $true }
    @{ Cell = 'E18'; Value = '  +4.53%  '; Numeric = $false }
    @{ Cell = 'D19'; Value = '0.000007866'; Numeric = $true }
    @{ Cell = 'E19'; Value = '  -0.24%  '; Numeric = $false }
    @{ Cell = 'D20'; Value = '13.26'; Numeric = $true }
    @{ Cell = 'E20'; Value = '  -0.08%  '; Numeric = $false }
    @{ Cell = 'D21'; Value = '0.9996'; Numeric = $true }
    @{ Cell = 'E21'; Value = '  -0.07%  '; Numeric = $false }
    @{ Cell = 'D22'; Value = '7.946'; Numeric = $true }
    @{ Cell = 'E22'; Value = '  +6.15%  '; Numeric = $false }
    @{ Cell = 'D23'; Value = '1.003'; Numeric = $true }
    @{ Cell = 'E23'; Value = '  +0.24%  '; Numeric = $false }
    @{ Cell = 'D24'; Value = '0.1591'; Numeric = $true }
    @{ Cell = 'E24'; Value = '  +13.03%  '; Numeric = $false }
    @{ Cell = 'D25'; Value = '163.92'; Numeric = $true }
    @{ Cell = 'E25'; Value = '  +0.45%  '; Numeric = $false }
    @{ Cell = 'D26'; Value = '9.001'; Numeric = $true }
    @{ Cell = 'E26'; Value = '  +1.34%  '; Numeric = $false }
    @{ Cell = 'D27'; Value = '18.31'; Numeric = $true }
    @{ Cell = 'E27'; Value = '  +1.29%  '; Numeric = $false }
    @{ Cell = 'E28'; Value = '  -3.90%  '; Numeric = $false }
    @{ Cell = 'D29'; Value = '1.499'; Numeric = $true }
    @{ Cell = 'E29'; Value = '  +1.67%  '; Numeric = $false }
    @{ Cell = 'D30'; Value = '4.375'; Numeric = $true }
    @{ Cell = 'E30'; Value = '  +0.32%  '; Numeric = $false }
    @{ Cell = 'D31'; Value = '4.119'; Numeric = $true }
    @{ Cell = 'E31'; Value = '  +2.43%  '; Numeric = $false }
    @{ Cell = 'D32'; Value = '0.05304'; Numeric = $true }
    @{ Cell = 'E32'; Value = '  +2.09%  '; Numeric = $false }
    @{ Cell = 'D33'; Value = '1.927'; Numeric = $true }
    @{ Cell = 'E33'; Value = '  +0.64%  '; Numeric = $false }
    @{ Cell = 'D34'; Value = '1.197'; Numeric = $true }
    @{ Cell = 'E34'; Value = '  +2.87%  '; Numeric = $false }
    @{ Cell = 'D35'; Value = '0.7247'; Numeric = $true }
    @{ Cell = 'E35'; Value = '  +1.24%  '; Numeric = $false }
    @{ Cell = 'D36'; Value = '2.678'; Numeric = $true }
    @{ Cell = 'E36'; Value = '  -0.37%  '; Numeric = $false }
    @{ Cell = 'D37'; Value = '0.01868'; Numeric = $true }
    @{ Cell = 'E37'; Value = '  +1.04%  '; Numeric = $false }
    @{ Cell = 'D38'; Value = '1.234.90'; Numeric = $false }
    @{ Cell = 'E38'; Value = '  +7.65%  '; Numeric = $false }
    @{ Cell = 'D39'; Value = '2.723'; Numeric = $true }
    @{ Cell = 'E39'; Value = '  +0.01%  '; Numeric = $false }
    @{ Cell = 'D40'; Value = '0.9072'; Numeric = $true }
    @{ Cell = 'E40'; Value = '  -2.63%  '; Numeric = $false }
    @{ Cell = 'D41'; Value = '74.01'; Numeric = $true }
    @{ Cell = 'E41'; Value = '  +5.13%  '; Numeric = $false }
    @{ Cell = 'D42'; Value = '6.164'; Numeric = $true }
    @{ Cell = 'E42'; Value = '  +2.44%  '; Numeric = $false }
    @{ Cell = 'D43'; Value = '1.001'; Numeric = $true }
    @{ Cell = 'E43'; Value = '  +0.03%  '; Numeric = $false }
    @{ Cell = 'D44'; Value = '103.17'; Numeric = $true }
    @{ Cell = 'E44'; Value = '  +0.39%  '; Numeric = $false }
    @{ Cell = 'D45'; Value = '0.5331'; Numeric = $true }
    @{ Cell = 'E45'; Value = '  +0.90%  '; Numeric = $false }
    @{ Cell = 'D46'; Value = '1.764'; Numeric = $true }
    @{ Cell = 'E46'; Value = '  +1.13%  '; Numeric = $false }
    @{ Cell = 'D47'; Value = '2.920'; Numeric = $true }
    @{ Cell = 'E47'; Value = '  +11.75%  '; Numeric = $false }
    @{ Cell = 'B48'; Value = 'EnergySwap'; Numeric = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Numeric = $false }
    @{ Cell = 'D48'; Value = '9.299'; Numeric = $true }
    @{ Cell = 'E48'; Value = '  +1.68%  '; Numeric = $false }
    @{ Cell = 'B49'; Value = 'TheSandbox'; Numeric = $false }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Numeric = $false }
    @{ Cell = 'D49'; Value = '0.4316'; Numeric = $true }
    @{ Cell = 'E49'; Value = '  +1.41%  '; Numeric = $false }
    @{ Cell = 'B50'; Value = 'Aptos'; Numeric = $false }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Numeric = $false }
    @{ Cell = 'D50'; Value = '7.058'; Numeric = $true }
    @{ Cell = 'E50'; Value = '  +1.48%  '; Numeric = $false }
    @{ Cell = 'B51'; Value = 'Frax'; Numeric = $false }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; Numeric = $false }
    @{ Cell = 'D51'; Value = '1.001'; Numeric = $true }
    @{ Cell = 'E51'; Value = '  +0.09%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $rng.Value = "'" + $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
